$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 2 blank columns before column D (shifts existing D:K to F:M)
$ws.Range("D:E").Insert(-4161)

# Step 2: propagate number formats from column F into new D:E columns for each data block
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: fix a handful of cells in F:M whose new values are not a pure shift of the old data
# (the source data refresh corrected/backfilled a few previously-blank "NA" cells)
$ws.Range("F8").Value = 0
$ws.Range("F18").Value = -9600
$ws.Range("F20").Value = 200
$ws.Range("F21").Value = -9400
$ws.Range("F32").Value = -200
$ws.Range("F91").Value = "NA"
$ws.Range("G91").Value = "NA"

# Step 4: populate the two new columns (D, E) with the refreshed quarterly data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = 0
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = 12000
$ws.Range("E12").Value = 13200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 16100
$ws.Range("E17").Value = 16800
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = -16800
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = 400
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = -16500
$ws.Range("D22").Value = 800
$ws.Range("E22").Value = 700
$ws.Range("D23").Value = -16500
$ws.Range("E23").Value = -17100
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -16500
$ws.Range("E26").Value = -17100
$ws.Range("D27").Value = -16500
$ws.Range("E27").Value = -17100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = -400
$ws.Range("D33").Value = -16500
$ws.Range("E33").Value = -17100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -16500
$ws.Range("E35").Value = -17100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 61300
$ws.Range("E41").Value = 16900
$ws.Range("D42").Value = 39100
$ws.Range("E42").Value = 48000
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 1500
$ws.Range("E45").Value = 1900
$ws.Range("D46").Value = 101800
$ws.Range("E46").Value = 66800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 200
$ws.Range("E48").Value = 200
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 400
$ws.Range("E52").Value = 400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 102400
$ws.Range("E54").Value = 67400
$ws.Range("D57").Value = 5800
$ws.Range("E57").Value = 6500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 6600
$ws.Range("D59").Value = 4200
$ws.Range("E59").Value = 1500
$ws.Range("D60").Value = 10000
$ws.Range("E60").Value = 14500
$ws.Range("D61").Value = 25600
$ws.Range("E61").Value = 18900
$ws.Range("D62").Value = 200
$ws.Range("E62").Value = 200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 35900
$ws.Range("E66").Value = 33600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -171200
$ws.Range("E72").Value = -154700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 66600
$ws.Range("E76").Value = 33800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -16500
$ws.Range("E81").Value = -17100
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -12500
$ws.Range("E89").Value = -13700
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 9100
$ws.Range("E94").Value = -7800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 47800
$ws.Range("E100").Value = 5100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 44300
$ws.Range("E102").Value = -16400
